$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the default (unstyled) cell style, used to
# restore style after forcing text entry via a quote-prefix so that
# numeric-looking values (e.g. "421.52") are stored as text, not numbers.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '67.257.25'
$ws.Range("E2").Value = '  +5.66%  '
$ws.Range("D3").Value = '3.713.75'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'421.52"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").Value = "'131.58"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  +1.82%  '
$ws.Range("D7").Value = '3.704.93'
$ws.Range("E7").Value = '  +6.99%  '
$ws.Range("D8").Value = "'0.645"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("D10").Value = "'0.766"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  -3.35%  '
$ws.Range("D11").Value = "'0.183"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  +12.18%  '
$ws.Range("D12").Value = "'0.0000400"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = '  +50.06%  '
$ws.Range("D13").Value = "'43.16"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").Value = "'10.26"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").Value = '4.304.43'
$ws.Range("E15").Value = '  +6.97%  '
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = "'20.80"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = '  +2.01%  '
$ws.Range("D18").Value = '3.707.33'
$ws.Range("E18").Value = '  +7.15%  '
$ws.Range("D19").Value = "'13.13"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  +5.71%  '
$ws.Range("E20").Value = '  +3.58%  '
$ws.Range("D21").Value = '67.340.75'
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("D22").Value = "'451.42"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = '  -4.02%  '
$ws.Range("D23").Value = "'15.86"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  +15.69%  '
$ws.Range("D24").Value = "'89.63"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").Value = "'3.17"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("D26").Value = "'38.14"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = '  +11.47%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").Value = "'10.18"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = "'3.33"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("E29").Value = '  +4.61%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = "'12.60"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = "'0.122"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = '  +8.85%  '
$ws.Range("D32").Value = "'2.79"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = '  +4.22%  '
$ws.Range("D33").Value = "'7.39"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").Value = "'41.92"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").Value = "'56.52"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").Value = "'0.0495"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '0.0₃0771'
$ws.Range("E39").Value = '  +18.97%  '
$ws.Range("D40").Value = "'3.18"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = '  +35.73%  '
$ws.Range("E41").Value = '  +3.76%  '
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = "'27.89"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = '  +26.67%  '
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").Value = "'3.43"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("D45").Value = "'148.17"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = "'2.11"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '  +5.36%  '
$ws.Range("D47").Value = "'2.92"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = '  -5.01%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'4.41"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  -2.16%  '
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").Value = "'0.308"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("E51").Value = '  +16.46%  '
